$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.363.20'
$ws.Range('E2').Value = '  -0.37%  '
$ws.Range('D3').Value = '3.883.49'
$ws.Range('E3').Value = '  -0.89%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = "'601.08"
$ws.Range('D6').Value = "'167.88"
$ws.Range('E6').Value = '  +1.31%  '
$ws.Range('D7').Value = '3.884.04'
$ws.Range('E7').Value = '  -0.94%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = "'0.528"
$ws.Range('E9').Value = '  -0.22%  '
$ws.Range('D10').Value = "'0.166"
$ws.Range('E10').Value = '  -0.04%  '
$ws.Range('D11').Value = "'6.48"
$ws.Range('E11').Value = '  +1.50%  '
$ws.Range('E12').Value = '  -0.19%  '
$ws.Range('E13').Value = '  -0.42%  '
$ws.Range('D14').Value = "'37.13"
$ws.Range('E14').Value = '  -0.41%  '
$ws.Range('D15').Value = '4.526.71'
$ws.Range('D16').Value = '3.889.47'
$ws.Range('E16').Value = '  -0.93%  '
$ws.Range('D17').Value = '68.306.50'
$ws.Range('E17').Value = '  -0.61%  '
$ws.Range('D18').Value = "'18.22"
$ws.Range('E18').Value = '  +6.18%  '
$ws.Range('D19').Value = "'7.42"
$ws.Range('E19').Value = '  -0.40%  '
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('E21').Value = '  -0.90%  '
$ws.Range('D22').Value = "'473.32"
$ws.Range('E22').Value = '  -2.65%  '
$ws.Range('D23').Value = "'0.735"
$ws.Range('E23').Value = '  +1.64%  '
$ws.Range('D24').Value = "'0.0000165"
$ws.Range('E24').Value = '  -3.17%  '
$ws.Range('D25').Value = "'83.72"
$ws.Range('E25').Value = '  -1.02%  '
$ws.Range('D26').Value = "'2.26"
$ws.Range('E26').Value = '  +1.38%  '
$ws.Range('D27').Value = "'12.21"
$ws.Range('E27').Value = '  +1.09%  '
$ws.Range('D28').Value = "'10.03"
$ws.Range('E28').Value = '  -0.83%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').Value = '  +1.13%  '
$ws.Range('D31').Value = '4.029.04'
$ws.Range('E31').Value = '  -1.00%  '
$ws.Range('D32').Value = "'7.94"
$ws.Range('E32').Value = '  +2.32%  '
$ws.Range('E33').Value = '  -2.68%  '
$ws.Range('D34').Value = "'31.48"
$ws.Range('E34').Value = '  -1.49%  '
$ws.Range('D35').Value = "'9.41"
$ws.Range('E35').Value = '  +1.99%  '
$ws.Range('D36').Value = '3.855.00'
$ws.Range('E36').Value = '  -0.36%  '
$ws.Range('E37').Value = '  -1.87%  '
$ws.Range('D38').Value = "'3.52"
$ws.Range('E38').Value = '  +11.00%  '
$ws.Range('D39').Value = "'1.03"
$ws.Range('E39').Value = '  -1.31%  '
$ws.Range('D40').Value = "'0.141"
$ws.Range('E40').Value = '  +2.23%  '
$ws.Range('D41').Value = "'5.93"
$ws.Range('E41').Value = '  +0.14%  '
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('D44').Value = "'434.97"
$ws.Range('E44').Value = '  +1.19%  '
$ws.Range('D45').Value = "'2.00"
$ws.Range('E45').Value = '  +0.71%  '
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('D47').Value = "'47.32"
$ws.Range('E47').Value = '  -2.34%  '
$ws.Range('D48').Value = "'8.60"
$ws.Range('E48').Value = '  +1.16%  '
$ws.Range('D49').Value = "'0.000286"
$ws.Range('E49').Value = '  +9.12%  '
$ws.Range('D50').Value = "'40.67"
$ws.Range('E50').Value = '  +4.29%  '
$ws.Range('D51').Value = "'143.75"
$ws.Range('E51').Value = '  +1.47%  '
